# Workbook is protected with a (known-blank/legacy) sheet password; unprotect
# so the cell writes below are not blocked by Excel's protected-sheet error.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure banner (A16).
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."

# Refresh Weight (D) and Percent Change (E) figures for rows 2-13.
$ws.Range("D2").Value = 0.03117872994871421
$ws.Range("E2").Value = -0.004045853000674393

$ws.Range("D3").Value = 0.02393167783955476
$ws.Range("E3").Value = 0.001391788448155795

$ws.Range("D4").Value = 0.05108145926576073
$ws.Range("E4").Value = 0.008018867924528328

$ws.Range("D5").Value = 0.1369536585288072
$ws.Range("E5").Value = 0.008149959250203676

$ws.Range("D6").Value = 0.02942002913372351
$ws.Range("E6").Value = 0.02852852852852861

$ws.Range("D7").Value = 0.1216077242045304
$ws.Range("E7").Value = -0.001603552485506365

$ws.Range("D8").Value = 0.1007925728729907
$ws.Range("E8").Value = 0.01036845028698385

$ws.Range("D9").Value = 0.02796227994241289
$ws.Range("E9").Value = 0.01658767772511838

$ws.Range("D10").Value = 0.1228373974204383
$ws.Range("E10").Value = 0.005769230769230749

$ws.Range("D11").Value = 0.2498406361937793
$ws.Range("E11").Value = -0.004131860235336404

$ws.Range("D12").Value = 0.104393834649288
$ws.Range("E12").Value = -0.006082493822467283

$ws.Range("E13").Value = 0.002627539175976024

# Re-apply protection to match the original workbook's protected-sheet state.
$ws.Protect()
